# Apply the "Add files via upload" revision to
# 工作进展/AI算法与应用-2024工作进展.xlsx
#
# Summary of the change:
#  - Sheet "004江润洲" (4th tab) becomes the active tab/selection (was
#    "009杨同学", the 9th/last tab).
#  - Sheet "004江润洲": row 4 grows new content in columns C/D/E, the row
#    height increases to fit it, and columns D/E get new/updated widths.
#  - Sheet "009杨同学": the previously selected whole-column selection
#    collapses back down to a single cell now that it's no longer the
#    active sheet.

$wb = $excel.ActiveWorkbook

$sheetJ = $wb.Worksheets.Item("004江润洲")
$sheetY = $wb.Worksheets.Item("009杨同学")

# ---------------------------------------------------------------------
# 1. New work-log content for 004江润洲, row 4 (columns C, D, E)
# ---------------------------------------------------------------------

$c4Text = @"
1. ASD多站点分类
（1）mobilenetv2更新一版结果，目前准确率0.882，ASD多站点分类sota准确率0.7
（2）分析网络权重，将权重映射到脑区找biomarker；
2. 自我学习
（1）双周6道题答案整理
"@.TrimEnd("`r","`n")

$d4Text = "1. 将权重映射到脑区找biomarker"

$e4Text = @"
1. ASD多站点分类
（1）通过目前的网络找到重要biomarker
（2）隐藏层特征可视化grad-CAM
2. 自我学习任务
（1）完成双周6道题
3. 其他任务
（1）筹备paper with code前沿与经典AI论文分享
（2）筹备AI算法与应用源码逐行解读
"@.TrimEnd("`r","`n")

# C4 and E4 reuse the wrapped/top-left alignment style already used by
# C2/C3/D2/D3/E2/E3 (cell style index 1) - copy it across before writing
# the values so no new style entries are created.
$sheetJ.Range("C2").Copy()
$sheetJ.Range("C4").PasteSpecial(-4122)   # xlPasteFormats
$sheetJ.Range("E4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$sheetJ.Range("C4").Value = $c4Text
$sheetJ.Range("D4").Value = $d4Text
$sheetJ.Range("E4").Value = $e4Text

# Row 4 grows taller to fit the new multi-line content.
$sheetJ.Rows.Item(4).RowHeight = 120.75

# New/updated column widths for D and E.
$sheetJ.Columns.Item(4).ColumnWidth = 34.714285714285715   # -> ~35.375 (D)
$sheetJ.Columns.Item(5).ColumnWidth = 49.857142857142854   # -> ~50.625 (E)

# ---------------------------------------------------------------------
# 2. Active tab / selection changes
# ---------------------------------------------------------------------

# 009杨同学 was the active sheet with the whole column C selected;
# collapse that selection down to the single top cell first.
$sheetY.Activate()
$sheetY.Range("C1").Select()

# 004江润洲 becomes the active sheet, selection moves to G4.
$sheetJ.Activate()
$sheetJ.Range("G4").Select()
